$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=0.1239850077872064;   C=0.9587971593447111; D=2.504666468680252; E=1.582613809076697; F=1.593449111284933; G=51}
    @{Row=3;  B=-0.07634650817340942; C=0.9644256602613136; D=2.465022428948611; E=1.570038989626885; F=1.584102668685676; G=50}
    @{Row=4;  B=0.1308097497450276;   C=0.9583096004403104; D=2.464810574400857; E=1.569971520251516; F=1.580725504613128; G=49}
    @{Row=5;  B=-0.02424481833496826; C=0.9746006404459205; D=2.576612619674879; E=1.605183048650489; F=1.621984539672506; G=48}
    @{Row=6;  B=0.1313683227808764;   C=1.001393086257909;  D=2.609607854279241; E=1.615428071527557; F=1.627484456937081; G=47}
    @{Row=7;  B=0.003180273963439154; C=0.9469229433371078; D=2.518249158100375; E=1.586899227455977; F=1.604431335894871; G=46}
    @{Row=8;  B=0.1362810442212024;   C=0.9540930716873162; D=2.589108079185797; E=1.60907056376835;  F=1.621405815584749; G=45}
    @{Row=9;  B=-0.04071899536153733; C=0.8760439952630932; D=2.432158382960957; E=1.559537874808097; F=1.577029996622433; G=44}
    @{Row=10; B=0.07889961332241763;  C=0.9613424106959066; D=2.595247230407774; E=1.610977104246915; F=1.628086442519692; G=43}
    @{Row=11; B=-0.01422600529044099; C=0.8653212001480697; D=2.498623644143271; E=1.580703528225097; F=1.599799477728101; G=42}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
}
